$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "62.594.82"
Set-TextValue "E2" "  -1.33%  "
Set-TextValue "D3" "3.008.09"
Set-TextValue "E3" "  -3.98%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "554.95"
Set-TextValue "E5" "  -1.17%  "
Set-TextValue "D6" "150.63"
Set-TextValue "E6" "  -7.18%  "
Set-TextValue "E7" "  +0.02%  "
Set-TextValue "E8" "  -3.39%  "
Set-TextValue "D9" "3.008.24"
Set-TextValue "E9" "  -3.84%  "
Set-TextValue "E10" "  -2.37%  "
Set-TextValue "D11" "6.34"
Set-TextValue "E11" "  -5.15%  "
Set-TextValue "D12" "0.365"
Set-TextValue "E12" "  -3.55%  "
Set-TextValue "D13" "3.527.89"
Set-TextValue "E13" "  -3.96%  "
Set-TextValue "E14" "  -3.37%  "
Set-TextValue "D15" "62.698.24"
Set-TextValue "E15" "  -1.16%  "
Set-TextValue "D16" "23.84"
Set-TextValue "E16" "  -4.14%  "
Set-TextValue "D17" "3.009.32"
Set-TextValue "E17" "  -3.54%  "
Set-TextValue "D18" "0.0000148"
Set-TextValue "E18" "  -2.64%  "
Set-TextValue "D19" "395.14"
Set-TextValue "E19" "  -1.45%  "
Set-TextValue "E20" "  -2.19%  "
Set-TextValue "D21" "11.87"
Set-TextValue "E21" "  -4.65%  "
Set-TextValue "E22" "  -5.72%  "
Set-TextValue "D23" "0.999"
Set-TextValue "E23" "  -0.03%  "
Set-TextValue "D24" "64.93"
Set-TextValue "E24" "  -3.47%  "
Set-TextValue "E25" "  -2.75%  "
Set-TextValue "D26" "0.186"
Set-TextValue "E26" "  -7.05%  "
Set-TextValue "D27" "0.0₃0964"
Set-TextValue "E27" "  -4.12%  "
Set-TextValue "D28" "8.58"
Set-TextValue "E28" "  -1.24%  "
Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  +0.05%  "
Set-TextValue "E30" "  +0.04%  "
Set-TextValue "E31" "  -2.53%  "
Set-TextValue "D32" "20.44"
Set-TextValue "E32" "  -2.06%  "
Set-TextValue "D33" "160.96"
Set-TextValue "E33" "  +5.46%  "
Set-TextValue "E34" "  -1.93%  "
Set-TextValue "D35" "6.01"
Set-TextValue "E35" "  -3.39%  "
Set-TextValue "E36" "  -2.00%  "
Set-TextValue "D37" "1.28"
Set-TextValue "E37" "  -2.80%  "
Set-TextValue "E38" "  -4.03%  "
Set-TextValue "D39" "2.468.06"
Set-TextValue "E39" "  -9.73%  "
Set-TextValue "D40" "37.52"
Set-TextValue "E40" "  -2.74%  "
Set-TextValue "D41" "22.47"
Set-TextValue "E41" "  -3.35%  "
Set-TextValue "D42" "3.90"
Set-TextValue "E42" "  -3.77%  "
Set-TextValue "E43" "  -4.27%  "
Set-TextValue "E44" "  -3.28%  "
Set-TextValue "E45" "  -0.14%  "
Set-TextValue "E46" "  -3.75%  "
Set-TextValue "D47" "4.96"
Set-TextValue "E47" "  -7.62%  "
Set-TextValue "D48" "19.70"
Set-TextValue "E48" "  -5.03%  "
Set-TextValue "B49" "Stellar"
Set-TextValue "C49" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D49" "0.0945"
Set-TextValue "E49" "  -2.88%  "
Set-TextValue "B50" "WhiteBITCoin"
Set-TextValue "C50" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D50" "10.49"
Set-TextValue "E50" "  +0.20%  "
Set-TextValue "D51" "261.96"
Set-TextValue "E51" "  -6.59%  "
